$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.000.98"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.679.14"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'215.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "'0.519"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D10").Value = "'20.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "1.915.39"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "1.670.31"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").Value = "'4.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "'65.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "'8.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.10%  "
$ws.Range("D18").Value = "27.017.91"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "'236.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("E24").Value = "  -2.80%  "
$ws.Range("D25").Value = "'145.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("D27").Value = "'16.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D28").Value = "'0.112"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D33").Value = "1.483.53"
$ws.Range("E33").Value = "  +2.35%  "
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("E35").Value = "  +5.11%  "
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").Value = "'0.584"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.67%  "
$ws.Range("E38").Value = "  +2.49%  "
$ws.Range("D39").Value = "'0.906"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.43%  "
$ws.Range("D40").Value = "'5.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.86%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("D44").Value = "'67.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("D45").Value = "1.820.41"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "'0.783"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "'90.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").Value = "'1.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("D51").Value = "'0.0509"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "
